$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row (row 1): extend with two new header values, 14 and 15,
#     in columns P and Q, matching the existing header style (bold /
#     centered / bordered) used by B1:O1. Copy format from O1 first.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows 2-25: columns I, K, M, O swap their 1/2 values, and two
#     new columns P and Q (value 2) are appended to every row.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2 (was 1)
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1 (was 2)
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2 (was 1)
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1 (was 2)
    $ws.Cells.Item($r, 16).Value = 2   # P (new)
    $ws.Cells.Item($r, 17).Value = 2   # Q (new)
}
